$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths.
# Note: the runtime's Range.ColumnWidth setter quantizes internally to whole
# "points" (pt = round(ColumnWidth*6); stored width = pt/6 + 5/6), so we pick
# the ColumnWidth value whose quantized result lands closest to the exact
# target OOXML column widths (13.7109375 and 14.42578125).
$ws.Columns.Item(1).ColumnWidth = 12.833333333333334
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666

# Update cell values (rows 1-3, row 4 unchanged)
$ws.Range("A1").Value = -0.077979173099918461
$ws.Range("B1").Value = 0.077979172742905112

$ws.Range("A2").Value = 0.068331130614467439
$ws.Range("B2").Value = -0.068331130973545276

$ws.Range("A3").Value = 0.025781416293584228
$ws.Range("B3").Value = -0.025781416680986914
